$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text edits (Volume number, report date range) ---
$ws.Range("A8").Value = "Volume 31   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/22/2024  Through  4/28/2024"

# --- Style-changing cells ---
# Style 14 donor (text/dash cells): C14
# Style 15 donor (plain count cells): I14
# Style 16 donor (percent-change cells): K14

$ws.Range("I14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F22").PasteSpecial(-4122)

$ws.Range("I14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 2

$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = 50

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D33").PasteSpecial(-4122)

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E33").PasteSpecial(-4122)

# --- Simple numeric value updates (style unchanged) ---
$ws.Range("N14").Value = -96
$ws.Range("D15").Value = 6
$ws.Range("E15").Value = -83.333333333333
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 9
$ws.Range("H15").Value = -77.777777777777
$ws.Range("I15").Value = 9
$ws.Range("J15").Value = 18
$ws.Range("K15").Value = -50
$ws.Range("L15").Value = -57.142857142857
$ws.Range("M15").Value = -10
$ws.Range("N15").Value = -62.5
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = 44.444444444444
$ws.Range("F16").Value = 41
$ws.Range("H16").Value = 24.242424242424
$ws.Range("I16").Value = 158
$ws.Range("J16").Value = 173
$ws.Range("K16").Value = -8.670520231213
$ws.Range("L16").Value = -16.402116402116
$ws.Range("M16").Value = 8.965517241379
$ws.Range("N16").Value = -75.914634146341
$ws.Range("C17").Value = 16
$ws.Range("D17").Value = 21
$ws.Range("E17").Value = -23.809523809523
$ws.Range("F17").Value = 53
$ws.Range("G17").Value = 70
$ws.Range("H17").Value = -24.285714285714
$ws.Range("I17").Value = 211
$ws.Range("J17").Value = 262
$ws.Range("K17").Value = -19.465648854961
$ws.Range("L17").Value = -0.471698113207
$ws.Range("M17").Value = 31.875
$ws.Range("N17").Value = -19.465648854961
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 25
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 93
$ws.Range("J18").Value = 125
$ws.Range("K18").Value = -25.6
$ws.Range("L18").Value = 3.333333333333
$ws.Range("M18").Value = -25
$ws.Range("N18").Value = -83.422459893048
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 27
$ws.Range("E19").Value = -55.555555555555
$ws.Range("F19").Value = 57
$ws.Range("G19").Value = 89
$ws.Range("H19").Value = -35.955056179775
$ws.Range("I19").Value = 277
$ws.Range("J19").Value = 322
$ws.Range("K19").Value = -13.975155279503
$ws.Range("L19").Value = -18.289085545722
$ws.Range("M19").Value = 68.90243902439
$ws.Range("N19").Value = 13.524590163934
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 32
$ws.Range("G20").Value = 57
$ws.Range("H20").Value = -43.859649122807
$ws.Range("I20").Value = 146
$ws.Range("J20").Value = 215
$ws.Range("K20").Value = -32.093023255814
$ws.Range("L20").Value = -22.340425531914
$ws.Range("M20").Value = 121.212121212121
$ws.Range("N20").Value = -76.751592356687
$ws.Range("C21").Value = 56
$ws.Range("D21").Value = 80
$ws.Range("E21").Value = -30
$ws.Range("F21").Value = 200
$ws.Range("G21").Value = 283
$ws.Range("H21").Value = -29.328621908127
$ws.Range("I21").Value = 895
$ws.Range("J21").Value = 1118
$ws.Range("K21").Value = -19.94633273703
$ws.Range("L21").Value = -14.189837008629
$ws.Range("M21").Value = 33.38301043219
$ws.Range("N21").Value = -62.708333333333
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -100
$ws.Range("C23").Value = 8
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 24
$ws.Range("G23").Value = 14
$ws.Range("H23").Value = 71.428571428571
$ws.Range("I23").Value = 84
$ws.Range("J23").Value = 92
$ws.Range("K23").Value = -8.695652173913
$ws.Range("L23").Value = -21.495327102803
$ws.Range("M23").Value = 23.529411764705
$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = 29.032258064516
$ws.Range("F24").Value = 111
$ws.Range("G24").Value = 128
$ws.Range("H24").Value = -13.28125
$ws.Range("I24").Value = 556
$ws.Range("J24").Value = 640
$ws.Range("K24").Value = -13.125
$ws.Range("L24").Value = -6.554621848739
$ws.Range("M24").Value = 24.663677130044
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = -15.384615384615
$ws.Range("F25").Value = 39
$ws.Range("G25").Value = 54
$ws.Range("H25").Value = -27.777777777777
$ws.Range("I25").Value = 196
$ws.Range("J25").Value = 279
$ws.Range("K25").Value = -29.749103942652
$ws.Range("L25").Value = -44.632768361581
$ws.Range("C26").Value = 19
$ws.Range("D26").Value = 22
$ws.Range("E26").Value = -13.636363636363
$ws.Range("F26").Value = 80
$ws.Range("G26").Value = 81
$ws.Range("H26").Value = -1.234567901234
$ws.Range("I26").Value = 370
$ws.Range("J26").Value = 336
$ws.Range("K26").Value = 10.119047619047
$ws.Range("L26").Value = 5.413105413105
$ws.Range("M26").Value = -22.105263157894
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 7
$ws.Range("E27").Value = -85.714285714285
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 12
$ws.Range("H27").Value = -58.333333333333
$ws.Range("I27").Value = 17
$ws.Range("J27").Value = 25
$ws.Range("K27").Value = -32
$ws.Range("L27").Value = -43.333333333333
$ws.Range("C28").Value = 3
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 120
$ws.Range("I28").Value = 40
$ws.Range("J28").Value = 35
$ws.Range("K28").Value = 14.285714285714
$ws.Range("L28").Value = 100
$ws.Range("N29").Value = -80
$ws.Range("N30").Value = -80.95238095238

$excel.CutCopyMode = $false
